# Ticket_Analysis.xlsx update
#  1. Reorder sheet tabs: "feedbacks" moves in front of "ticket"
#     (new order: feedbacks, ticket, Results)
#  2. "feedbacks" sheet gains a lookup formula in column B (ticket_created_at)
#     pulled from the "ticket" sheet, formatted as a custom date/time format
#  3. "Results" sheet gains SUMPRODUCT array formulas for same_day_count /
#     same_hour_count in columns B and C
#  4. Active-tab / selection bookkeeping: "feedbacks" becomes the active tab,
#     "Results" is no longer active

$wb = $excel.ActiveWorkbook

# --- 1. Move "feedbacks" to be the first sheet tab -------------------------
$wb.Worksheets.Item("feedbacks").Move($wb.Worksheets.Item(1))

# --- 2. Update the "Results" sheet's selection (it stops being active) -----
$wsResults = $wb.Worksheets.Item("Results")
$wsResults.Range("F4").Select()

# --- 3. Fill in the new lookup formulas on "feedbacks" ----------------------
$wsFeedbacks = $wb.Worksheets.Item("feedbacks")

$wsFeedbacks.Range("B2").Formula = '=IFERROR(INDEX(ticket!$B$2:$B$1000, MATCH(A2, ticket!$E$2:$E$1000, 0)), "")'
$wsFeedbacks.Range("B3").Formula = '=IFERROR(INDEX(ticket!$B$2:$B$1000, MATCH(A3, ticket!$E$2:$E$1000, 0)), "")'

# Custom number format dd/mm/yyyy hh:mm:ss applied to the new formula cells
$wsFeedbacks.Range("B2:B3").NumberFormat = 'dd/mm/yyyy\ hh:mm:ss'

# Column B widened to fit the new date/time values
$wsFeedbacks.Columns.Item(2).ColumnWidth = 17.3

# --- 4. "feedbacks" becomes the active sheet/tab, with its own selection ---
$wsFeedbacks.Activate()
$wsFeedbacks.Range("F11").Select()

# --- 5. Add the SUMPRODUCT array formulas on "Results" ----------------------
$wsResults.Range("B2").FormulaArray = '=SUMPRODUCT((ticket!$D$2:$D$1000 = $A2) * (INT(ticket!$B$2:$B$1000) = INT(ticket!$C$2:$C$1000)))'
$wsResults.Range("C2").FormulaArray = '=SUMPRODUCT((ticket!$D$2:$D$1000 = $A2) * (INT(ticket!$B$2:$B$1000) = INT(ticket!$C$2:$C$1000)) * (HOUR(ticket!$B$2:$B$1000) = HOUR(ticket!$C$2:$C$1000)))'
$wsResults.Range("B3").FormulaArray = '=SUMPRODUCT((ticket!$D$2:$D$1000 = $A3) * (INT(ticket!$B$2:$B$1000) = INT(ticket!$C$2:$C$1000)))'
$wsResults.Range("C3").FormulaArray = '=SUMPRODUCT((ticket!$D$2:$D$1000 = $A3) * (INT(ticket!$B$2:$B$1000) = INT(ticket!$C$2:$C$1000)) * (HOUR(ticket!$B$2:$B$1000) = HOUR(ticket!$C$2:$C$1000)))'
